# feat: added delete student
# Insert new columns around the existing email/status columns so the sheet
# grows from (email, status) in A:B to the full student-record layout in A:J,
# then fill in the new header cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room for two new leading columns (studentID, username) before the
#    existing "email" column (currently A).
$ws.Range("A:B").Insert()

# 2) Make room for five new columns (contactNo..studentPGMarks) between the
#    "email" column (now C) and the "status" column (now D).
$ws.Range("D:H").Insert()

# 3) Make room for one trailing column (studentDescription) after "status"
#    (now I).
$ws.Range("J:J").Insert()

# Whole-column inserts drag the neighbouring column's formatting along for
# the ride (Excel's own "format from left" behaviour); the authored sheet
# only carries the hyperlink style on the email column itself, so strip the
# borrowed style back off the newly inserted, still-empty cells.
$ws.Range("D2:H6").Clear()

# New header values for the inserted columns.
$ws.Range("A1").Value = "studentID"
$ws.Range("B1").Value = "username"
$ws.Range("D1").Value = "contactNo"
$ws.Range("E1").Value = "tenthMarks"
$ws.Range("F1").Value = "twelthMarks"
$ws.Range("G1").Value = "studentUGMarks"
$ws.Range("H1").Value = "studentPGMarks"
$ws.Range("J1").Value = "studentDescription"

# Dimension grows to A1:J6 automatically; update the active selection like
# the authored workbook (cursor left on H9 after the edits).
$ws.Range("H9").Select()
